# Apply the "update new orleans xlsx files" edit:
#  1. hotel_info gains a new "State" column (value "Louisiana") inserted
#     right after "Hotel_Name" and before "City".
#  2. The "review_info" sheet tab is moved in front of "hotel_info"
#     (review_info becomes the first/active sheet tab).

$wb = $excel.ActiveWorkbook

$wsHotel  = $wb.Worksheets.Item("hotel_info")
$wsReview = $wb.Worksheets.Item("review_info")

# --- 1. Insert the new State column into hotel_info ---------------------
# Hotel_Name is column B, City is column C -> insert a new column C so the
# new column sits between them, then fill in header + value.
$wsHotel.Columns.Item(3).Insert()
$wsHotel.Range("C1").Value = "State"
$wsHotel.Range("C2").Value = "Louisiana"

# --- 2. Reorder sheet tabs: review_info before hotel_info ---------------
$wsReview.Move($wsHotel)
